# Update workbook to reflect data through 2021-10-08 (adds data for 2021-10-16... per commit msg,
# but the concrete OOXML diff we must reproduce touches the header labels and a handful of cell
# values in the carjacking-by-neighborhood-by-month sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the "through" date label (shared string + header cell).
$ws.Name = "Through 2021-10-08"
$ws.Range("B1").Value = "October 2021 (through October 08)"

# Row 2 - Garfield Park
$ws.Range("AF2").Value = 1
$ws.Range("AZ2").Value = 3

# Row 3 - Austin
$ws.Range("B3").Value = 3
$ws.Range("V3").Value = 2

# Row 4 - North Lawndale
$ws.Range("B4").Value = 4

# Row 6 - Auburn Gresham
$ws.Range("AZ6").Value = 1

# Row 9 - Grand Crossing
$ws.Range("AZ9").Value = 1

# Row 14 - United Center
$ws.Range("AZ14").Value = 1

# Row 18 - Little Village
$ws.Range("AZ18").Value = 1

# Row 28 - West Loop
$ws.Range("AF28").Value = 1

# Row 41 - Uptown
$ws.Range("AF41").Value = 2

# Row 48 - Washington Heights
$ws.Range("AP48").Value = 1
